$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename (row 1): Spanish labels -> snake_case machine names ---
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# --- Title-case the Spanish connector words ("de"/"del"/"la"/"el"/"los"/"las"/"y") ---
# --- inside state/municipality names throughout the data body ---
$ws.Range("B8").Value = 'Pabellón De Arteaga'
$ws.Range("B9").Value = 'Rincón De Romos'
$ws.Range("B10").Value = 'San Francisco De Los Romo'
$ws.Range("B11").Value = 'San José De Gracia'
$ws.Range("B36").Value = 'Chiapa De Corzo'
$ws.Range("B38").Value = 'Comitán De Domínguez'
$ws.Range("B54").Value = 'Mazapa De Madero'
$ws.Range("B58").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B88").Value = 'Coyame Del Sotol'
$ws.Range("B98").Value = 'Guadalupe Y Calvo'
$ws.Range("B100").Value = 'Hidalgo Del Parral'
$ws.Range("B121").Value = 'San Francisco De Borja'
$ws.Range("B122").Value = 'San Francisco De Conchos'
$ws.Range("B123").Value = 'San Francisco Del Oro'
$ws.Range("B130").Value = 'Valle De Zaragoza'
$ws.Range("B151").Value = 'San Juan De Sabinas'
$ws.Range("B164").Value = 'Villa De Álvarez'
$ws.Range("A166").Value = 'Ciudad De México'
$ws.Range("B170").Value = 'Cuajimalpa De Morelos'
$ws.Range("B184").Value = 'Coneto De Comonfort'
$ws.Range("B198").Value = 'Nombre De Dios'
$ws.Range("B201").Value = 'Pánuco De Coronado'
$ws.Range("B208").Value = 'San Juan De Guadalupe'
$ws.Range("B209").Value = 'San Juan Del Río'
$ws.Range("B210").Value = 'San Luis Del Cordero'
$ws.Range("B211").Value = 'San Pedro Del Gallo'
$ws.Range("A221").Value = 'Estado De México'
$ws.Range("B221").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B226").Value = 'Atizapán De Zaragoza'
$ws.Range("B232").Value = 'Coacalco De Berriozábal'
$ws.Range("B237").Value = 'Ecatepec De Morelos'
$ws.Range("B240").Value = 'Ixtapan De La Sal'
$ws.Range("B248").Value = 'Naucalpan De Juárez'
$ws.Range("B254").Value = 'San Antonio La Isla'
$ws.Range("B255").Value = 'San Felipe Del Progreso'
$ws.Range("B257").Value = 'Soyaniquilpan De Juárez'
$ws.Range("B263").Value = 'Tenango Del Valle'
$ws.Range("B266").Value = 'Tlalnepantla De Baz'
$ws.Range("B271").Value = 'Valle De Bravo'
$ws.Range("B272").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B273").Value = 'Villa De Allende'
$ws.Range("B274").Value = 'Villa Del Carbón'
$ws.Range("B282").Value = 'San Miguel De Allende'
$ws.Range("B283").Value = 'Apaseo El Alto'
$ws.Range("B284").Value = 'Apaseo El Grande'
$ws.Range("B292").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B296").Value = 'Jaral Del Progreso'
$ws.Range("B304").Value = 'Purísima Del Rincón'
$ws.Range("B307").Value = 'San Diego De La Unión'
$ws.Range("B309").Value = 'San Francisco Del Rincón'
$ws.Range("B311").Value = 'San Luis De La Paz'
$ws.Range("B312").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B313").Value = 'Silao De La Victoria'
$ws.Range("B318").Value = 'Valle De Santiago'
$ws.Range("B324").Value = 'Acapulco De Juárez'
$ws.Range("B326").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B327").Value = 'Alcozauca De Guerrero'
$ws.Range("B332").Value = 'Atoyac De Álvarez'
$ws.Range("B333").Value = 'Ayutla De Los Libres'
$ws.Range("B336").Value = 'Chilapa De Álvarez'
$ws.Range("B337").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B338").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B341").Value = 'Coyuca De Benítez'
$ws.Range("B342").Value = 'Coyuca De Catalán'
$ws.Range("B345").Value = 'Cuetzala Del Progreso'
$ws.Range("B346").Value = 'Cutzamala De Pinzón'
$ws.Range("B351").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B352").Value = 'Iguala De La Independencia'
$ws.Range("B353").Value = 'Zihuatanejo De Azueta'
$ws.Range("B355").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B358").Value = 'Mártir De Cuilapan'
$ws.Range("B367").Value = 'Taxco De Alarcón'
$ws.Range("B369").Value = 'Técpan De Galeana'
$ws.Range("B371").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B373").Value = 'Tixtla De Guerrero'
$ws.Range("B374").Value = 'Tlapa De Comonfort'
$ws.Range("B387").Value = 'Cuautepec De Hinojosa'
$ws.Range("B391").Value = 'Huasca De Ocampo'
$ws.Range("B393").Value = 'Huejutla De Reyes'
$ws.Range("B396").Value = 'Jacala De Ledezma'
$ws.Range("B400").Value = 'Mixquiahuala De Juárez'
$ws.Range("B401").Value = 'Molango De Escamilla'
$ws.Range("B403").Value = 'Nopala De Villagrán'
$ws.Range("B404").Value = 'Pachuca De Soto'
$ws.Range("B406").Value = 'Progreso De Obregón'
$ws.Range("B409").Value = 'Santiago De Anaya'
$ws.Range("B410").Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range("B413").Value = 'Tenango De Doria'
$ws.Range("B415").Value = 'Tepehuacán De Guerrero'
$ws.Range("B416").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B418").Value = 'Tezontepec De Aldama'
$ws.Range("B422").Value = 'Tula De Allende'
$ws.Range("B423").Value = 'Tulancingo De Bravo'
$ws.Range("B425").Value = 'Zacualtipán De Ángeles'
$ws.Range("B433").Value = 'Atemajac De Brizuela'
$ws.Range("B434").Value = 'Atotonilco El Alto'
$ws.Range("B435").Value = 'Autlán De Navarro'
$ws.Range("B440").Value = 'Cañadas De Obregón'
$ws.Range("B445").Value = 'Concepción De Buenos Aires'
$ws.Range("B451").Value = 'Encarnación De Díaz'
$ws.Range("B455").Value = 'Huejuquilla El Alto'
$ws.Range("B461").Value = 'Lagos De Moreno'
$ws.Range("B465").Value = 'Ojuelos De Jalisco'
$ws.Range("B469").Value = 'San Cristóbal De La Barranca'
$ws.Range("B471").Value = 'San Juan De Los Lagos'
$ws.Range("B472").Value = 'San Juanito De Escobedo'
$ws.Range("B474").Value = 'San Martín De Bolaños'
$ws.Range("B476").Value = 'Santa María De Los Ángeles'
$ws.Range("B478").Value = 'Tamazula De Gordiano'
$ws.Range("B481").Value = 'Tepatitlán De Morelos'
$ws.Range("B491").Value = 'Unión De San Antonio'
$ws.Range("B492").Value = 'Unión De Tula'
$ws.Range("B493").Value = 'Valle De Juárez'
$ws.Range("B497").Value = 'Yahualica De González Gallo'
$ws.Range("B498").Value = 'Zacoalco De Torres'
$ws.Range("B501").Value = 'Zapotlán Del Rey'
$ws.Range("B502").Value = 'Zapotlán El Grande'
$ws.Range("B520").Value = 'Coalcomán De Vázquez Pallares'
$ws.Range("B562").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B582").Value = 'Coatlán Del Río'
$ws.Range("B589").Value = 'Jonacatepec De Leandro Valle'
$ws.Range("B592").Value = 'Puente De Ixtla'
$ws.Range("B597").Value = 'Tetela Del Volcán'
$ws.Range("B611").Value = 'Santa María Del Oro'
$ws.Range("B622").Value = 'Lampazos De Naranjo'
$ws.Range("B624").Value = 'Mier Y Noriega'
$ws.Range("B629").Value = 'San Nicolás De Los Garza'
$ws.Range("B634").Value = 'Ayoquezco De Aldama'
$ws.Range("B639").Value = 'El Barrio De La Soledad'
$ws.Range("B640").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B641").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B642").Value = 'Ixtlán De Juárez'
$ws.Range("B643").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B647").Value = 'Mariscala De Juárez'
$ws.Range("B649").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B650").Value = 'Oaxaca De Juárez'
$ws.Range("B651").Value = 'Ocotlán De Morelos'
$ws.Range("B652").Value = 'Putla Villa De Guerrero'
$ws.Range("B657").Value = 'San Dionisio Del Mar'
$ws.Range("B661").Value = 'San José Del Progreso'
$ws.Range("B674").Value = 'San Miguel Del Puerto'
$ws.Range("B688").Value = 'Santa María Del Tule'
$ws.Range("B708").Value = 'Santo Domingo De Morelos'
$ws.Range("B715").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B716").Value = 'Teotitlán De Flores Magón'
$ws.Range("B717").Value = 'Tlacolula De Matamoros'
$ws.Range("B719").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B720").Value = 'Zimatlán De Álvarez'
$ws.Range("B735").Value = 'Cuetzalan Del Progreso'
$ws.Range("B741").Value = 'Ixcamilpa De Guerrero'
$ws.Range("B743").Value = 'Izúcar De Matamoros'
$ws.Range("B760").Value = 'San Salvador El Verde'
$ws.Range("B767").Value = 'Tepexi De Rodríguez'
$ws.Range("B768").Value = 'Tetela De Ocampo'
$ws.Range("B776").Value = 'Tuzamapan De Galeana'
$ws.Range("B785").Value = 'Amealco De Bonfil'
$ws.Range("B787").Value = 'Cadereyta De Montes'
$ws.Range("B790").Value = 'Jalpan De Serra'
$ws.Range("B791").Value = 'Landa De Matamoros'
$ws.Range("B792").Value = 'Pinal De Amoles'
$ws.Range("B794").Value = 'San Juan Del Río'
$ws.Range("B804").Value = 'Armadillo De Los Infante'
$ws.Range("B809").Value = 'Cerro De San Pedro'
$ws.Range("B811").Value = 'Ciudad Del Maíz'
$ws.Range("B820").Value = 'Mexquitic De Carmona'
$ws.Range("B825").Value = 'San Ciro De Acosta'
$ws.Range("B831").Value = 'Santa María Del Río'
$ws.Range("B833").Value = 'Soledad De Graciano Sánchez'
$ws.Range("B840").Value = 'Villa De Arista'
$ws.Range("B841").Value = 'Villa De Arriaga'
$ws.Range("B842").Value = 'Villa De Guadalupe'
$ws.Range("B843").Value = 'Villa De Ramos'
$ws.Range("B844").Value = 'Villa De Reyes'
$ws.Range("B914").Value = 'Soto La Marina'
$ws.Range("B924").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B925").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B926").Value = 'San Pablo Del Monte'
$ws.Range("B927").Value = 'Tepetitla De Lardizábal'
$ws.Range("B936").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B943").Value = 'Boca Del Río'
$ws.Range("B945").Value = 'Castillo De Teayo'
$ws.Range("B946").Value = 'Cazones De Herrera'
$ws.Range("B956").Value = 'Cosamaloapan De Carpio'
$ws.Range("B968").Value = 'Hueyapan De Ocampo'
$ws.Range("B969").Value = 'Ignacio De La Llave'
$ws.Range("B971").Value = 'Ixhuacán De Los Reyes'
$ws.Range("B977").Value = 'Juchique De Ferrer'
$ws.Range("B980").Value = 'Lerdo De Tejada'
$ws.Range("B982").Value = 'Martínez De La Torre'
$ws.Range("B991").Value = 'Ozuluama De Mascareñas'
$ws.Range("B997").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1003").Value = 'Sayula De Alemán'
$ws.Range("B1004").Value = 'Soledad De Doblado'
$ws.Range("B1006").Value = 'Tatahuicapan De Juárez'
$ws.Range("B1021").Value = 'Tlacotepec De Mejía'
$ws.Range("B1026").Value = 'Vega De Alatorre'
$ws.Range("B1032").Value = 'Zozocolco De Hidalgo'
$ws.Range("B1043").Value = 'Cañitas De Felipe Pescador'
$ws.Range("B1045").Value = 'Concepción Del Oro'
$ws.Range("B1066").Value = 'Nochistlán De Mejía'
$ws.Range("B1077").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1080").Value = 'Villa De Cos'

# --- Drop the trailing metadata/footer rows (old rows 1088-1092) and shrink the ---
# --- used range back down to A1:D1086 ---
$ws.Range("A1088:A1092").EntireRow.Delete()

